# Revert "Migration to Automation-Org/TestCases-maintenance/WIP-RMA TestCases"
# Removes the extra Labor Booking User / SiteID / Location ID / Location Number
# columns (H:K) that the migration had introduced, restores column B's width,
# flips the Background Processing flag back, and restores the prior selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the migrated-in columns H:K entirely (headers + data + shared strings
# for them are cleaned up automatically).
$ws.Range("H1:K3").EntireColumn.Delete()

# Column B goes back to its wider, pre-migration width.
$ws.Columns("B").ColumnWidth = 37.1666666666667

# Swap the Background Processing boolean flag back on row 2 / off on row 3.
$ws.Range("G2").Value = $true
$ws.Range("G3").Value = $false

# Restore the previously-active selection.
$ws.Range("E12").Select() | Out-Null
